$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.07702833333333334
$ws.Range("H2").Value = 0.231085
$ws.Range("I2").Value = 0.02259036512642383
$ws.Range("J2").Value = 0.02259036512642383
$ws.Range("M2").Value = 0.6100786666666667
$ws.Range("N2").Value = 1.830236
$ws.Range("O2").Value = 0.06402955811028149
$ws.Range("P2").Value = 0.06402955811028149
$ws.Range("Q2").Value = 0.04699334289555556
$ws.Range("R2").Value = 0.42294008606
$ws.Range("S2").Value = 0.001446451096594831
$ws.Range("T2").Value = 0.001446451096594831
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.07702833333333334
$ws.Range("H3").Value = 0.231085
$ws.Range("I3").Value = 0.02259036512642383
$ws.Range("J3").Value = 0.02259036512642383
$ws.Range("M3").Value = 7.236132333333333
$ws.Range("O3").Value = 0.7594534623909487
$ws.Range("P3").Value = 0.7594534623909487
$ws.Range("Q3").Value = 0.5573872134161111
$ws.Range("R3").Value = 5.016484920745
$ws.Range("S3").Value = 0.01715633101193832
$ws.Range("T3").Value = 0.01715633101193832
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.07702833333333334
$ws.Range("H4").Value = 0.231085
$ws.Range("I4").Value = 0.02259036512642383
$ws.Range("J4").Value = 0.02259036512642383
$ws.Range("M4").Value = 1.681867666666667
$ws.Range("N4").Value = 5.045603
$ws.Range("O4").Value = 0.1765169794987699
$ws.Range("P4").Value = 0.1765169794987699
$ws.Range("Q4").Value = 0.1295514632505556
$ws.Range("R4").Value = 1.165963169255
$ws.Range("S4").Value = 0.003987583017890681
$ws.Range("T4").Value = 0.003987583017890682
$ws.Range("I5").Value = 0.9608869019286738
$ws.Range("J5").Value = 0.9608869019286738
$ws.Range("M5").Value = 0.6100786666666667
$ws.Range("N5").Value = 1.830236
$ws.Range("O5").Value = 0.06402955811028149
$ws.Range("P5").Value = 0.06402955811028149
$ws.Range("Q5").Value = 1.998873741680445
$ws.Range("R5").Value = 17.989863675124
$ws.Range("S5").Value = 0.06152516372445037
$ws.Range("T5").Value = 0.06152516372445037
$ws.Range("I6").Value = 0.9608869019286738
$ws.Range("J6").Value = 0.9608869019286738
$ws.Range("M6").Value = 7.236132333333333
$ws.Range("O6").Value = 0.7594534623909487
$ws.Range("P6").Value = 0.7594534623909487
$ws.Range("S6").Value = 0.7297488846358432
$ws.Range("T6").Value = 0.7297488846358432
$ws.Range("I7").Value = 0.9608869019286738
$ws.Range("J7").Value = 0.9608869019286738
$ws.Range("M7").Value = 1.681867666666667
$ws.Range("N7").Value = 5.045603
$ws.Range("O7").Value = 0.1765169794987699
$ws.Range("P7").Value = 0.1765169794987699
$ws.Range("Q7").Value = 5.510504299797445
$ws.Range("R7").Value = 49.594538698177
$ws.Range("S7").Value = 0.1696128535683802
$ws.Range("T7").Value = 0.1696128535683802
$ws.Range("G8").Value = 0.056339
$ws.Range("H8").Value = 0.169017
$ws.Range("I8").Value = 0.01652273294490242
$ws.Range("J8").Value = 0.01652273294490242
$ws.Range("M8").Value = 0.6100786666666667
$ws.Range("N8").Value = 1.830236
$ws.Range("O8").Value = 0.06402955811028149
$ws.Range("P8").Value = 0.06402955811028149
$ws.Range("Q8").Value = 0.03437122200133334
$ws.Range("R8").Value = 0.309340998012
$ws.Range("S8").Value = 0.001057943289236292
$ws.Range("T8").Value = 0.001057943289236292
$ws.Range("G9").Value = 0.056339
$ws.Range("H9").Value = 0.169017
$ws.Range("I9").Value = 0.01652273294490242
$ws.Range("J9").Value = 0.01652273294490242
$ws.Range("M9").Value = 7.236132333333333
$ws.Range("O9").Value = 0.7594534623909487
$ws.Range("P9").Value = 0.7594534623909487
$ws.Range("Q9").Value = 0.4076764595276666
$ws.Range("R9").Value = 3.669088135749
$ws.Range("S9").Value = 0.01254824674316714
$ws.Range("T9").Value = 0.01254824674316714
$ws.Range("G10").Value = 0.056339
$ws.Range("H10").Value = 0.169017
$ws.Range("I10").Value = 0.01652273294490242
$ws.Range("J10").Value = 0.01652273294490242
$ws.Range("M10").Value = 1.681867666666667
$ws.Range("N10").Value = 5.045603
$ws.Range("O10").Value = 0.1765169794987699
$ws.Range("P10").Value = 0.1765169794987699
$ws.Range("Q10").Value = 0.09475474247233333
$ws.Range("R10").Value = 0.852792682251
$ws.Range("S10").Value = 0.00291654291249899
$ws.Range("T10").Value = 0.002916542912498991
